$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 131, pushing the former rows 131-138 down to 132-139.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new weekly price record.
$ws.Cells.Item(131, 1).Value2  = 10
$ws.Cells.Item(131, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(131, 3).Value2  = "La Araucanía"
$ws.Cells.Item(131, 4).Value2  = 44714
$ws.Cells.Item(131, 5).Value2  = 9
$ws.Cells.Item(131, 6).Value2  = 100112031
$ws.Cells.Item(131, 7).Value2  = "Poroto verde"
$ws.Cells.Item(131, 8).Value2  = "Sin especificar"
$ws.Cells.Item(131, 9).Value2  = "Primera"
$ws.Cells.Item(131, 10).Value2 = 85
$ws.Cells.Item(131, 11).Value2 = 25000
$ws.Cells.Item(131, 12).Value2 = 25000
$ws.Cells.Item(131, 13).Value2 = 25000
$ws.Cells.Item(131, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(131, 15).Value2 = "Región del Maule"
$ws.Cells.Item(131, 16).Value2 = 1000
$ws.Cells.Item(131, 17).Value2 = 25
$ws.Cells.Item(131, 18).Value2 = "Hortaliza"
